$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$tail = $tr.Characters(24, 12)
$tail.Text = "using Python"
